$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries update: Ucrania overtakes Israel (row 27/28 swap names + new counts) ---
# Row 27 (was Israel) -> now Ucrania, with refreshed stats
$ws.Range("A27").Value = "Ucrania"
$ws.Range("B27").Value = 154335
$ws.Range("C27").Value = 2476
$ws.Range("D27").Value = 68893
$ws.Range("E27").Value = 82264
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 30
$ws.Range("H27").Value = 3178

# Row 28 (was Ucrania) -> now Israel, with refreshed stats
$ws.Range("A28").Value = "Israel"
$ws.Range("B28").Value = 153217
$ws.Range("C28").Value = 495
$ws.Range("D28").Value = 114624
$ws.Range("E28").Value = 37490
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 1103

# --- Provincias Spain (El Salvador row) refreshed stats ---
$ws.Range("D75").Value = 17885
$ws.Range("E75").Value = 8181
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 785

# --- Update the "last updated" timestamp string in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 08:18"
